# Updated symbol list on Mon Dec 26 22:40:11 UTC 2022 with GitHub Actions
#
# The nightly scraper re-ran: most rows just got a refreshed "Price" (column D)
# quote, but the coin that used to sit at rank #17 ("One") jumped up to rank #9,
# pushing WazirX / MandalaExchangeToken / LiechtensteinCryptoassetsExchange /
# BitrueCoin / BitMartToken / MCDex / BitForexToken / CoinExToken each down one
# spot (rows 10-18, columns B/C/E) with freshly scraped prices to match.
#
# Price cells are free-text numeric strings (e.g. "0.05910", "0.00005244")
# where trailing/leading zeros are significant, so each is explicitly
# formatted as Text before the value is written - otherwise Excel would
# coerce them to numbers and silently drop the significant zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Simple price refreshes (rows 2-9) ---
Set-PriceText "D2" "242.87"
Set-PriceText "D3" "23.11"
Set-PriceText "D4" "5.422"
Set-PriceText "D5" "0.05910"
Set-PriceText "D6" "3.448"
Set-PriceText "D7" "6.523"
Set-PriceText "D8" "0.8093"
Set-PriceText "D9" "0.9357"

# --- Rank reshuffle: "One" moves from #17 (row 18) up to #9 (row 10); ---
# --- WazirX..CoinExToken each shift down one row, with new prices.    ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-PriceText "D10" "0.0005925"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-PriceText "D11" "0.1426"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-PriceText "D12" "0.07420"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-PriceText "D13" "0.03249"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-PriceText "D14" "0.03094"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-PriceText "D15" "0.09359"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-PriceText "D16" "3.868"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-PriceText "D17" "0.001568"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-PriceText "D18" "0.04689"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- More simple price refreshes further down the table ---
Set-PriceText "D19" "0.005982"
Set-PriceText "D21" "0.004899"
Set-PriceText "D24" "2.141"
Set-PriceText "D40" "0.03957"
Set-PriceText "D41" "0.006586"
Set-PriceText "D42" "0.1073"
Set-PriceText "D43" "0.003003"
Set-PriceText "D44" "0.008774"
Set-PriceText "D45" "0.00005244"
Set-PriceText "D47" "0.6706"
Set-PriceText "D48" "0.002395"
